$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-9 from 2023-09-01 (45170) to 2023-09-05 (45174)
$ws.Range("C2:C9").Value = 45174
